# Extend age tables for lookup completeness
#
# The "IESD-FoIERbA" sheet (Fraction Retiring This Year) listed a lookup
# table of "Age" (row 1) vs "Fraction Retiring" (row 2) for ages 0-50
# (columns B:AZ). This extends that table out to age 210 (columns BA:HD)
# so that downstream XLOOKUP/MATCH formulas never run off the end of the
# table: ages beyond the last modeled age (50) are treated as 100%
# retired (value of 1), matching the saturated value already present in
# column AZ.

$wb = $excel.ActiveWorkbook

$wsAge = $wb.Worksheets.Item("IESD-FoIERbA")

for ($age = 51; $age -le 210; $age++) {
    $col = $age + 2   # column B (index 2) holds age 0, so age -> col = age + 2
    $wsAge.Cells.Item(1, $col).Value = $age
    $wsAge.Cells.Item(2, $col).Value = 1
}

# Row 2 (the header/description row for the retiring-fraction series) grew
# taller once the table was widened and re-wrapped in the newer Excel build.
$wsAge.Rows.Item(2).RowHeight = 45

# Reflect the author's final on-screen state: they ended up focused on the
# newly extended "IESD-FoIERbA" sheet (selecting a cell in the new range),
# after having been on "IESD-AAaWER" (where the selection was left at F20).
$wsOther = $wb.Worksheets.Item("IESD-AAaWER")
[void]$wsOther.Activate()
[void]$wsOther.Range("F20").Select()

[void]$wsAge.Activate()
[void]$wsAge.Range("AY3").Select()
try {
    $excel.ActiveWindow.ScrollColumn = 43
} catch {
}
